$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Updated rider counts (column C) and average values (column D)
$ws.Range("C2").Value = 219
$ws.Range("D2").Value = 209.12

$ws.Range("C3").Value = 288
$ws.Range("D3").Value = 208.59

$ws.Range("C4").Value = 241
$ws.Range("D4").Value = 230.38

$ws.Range("C5").Value = 258
$ws.Range("D5").Value = 239.82

$ws.Range("C6").Value = 98
$ws.Range("D6").Value = 118.83

$ws.Range("C7").Value = 93
$ws.Range("D7").Value = 101.88
